$wb = $excel.ActiveWorkbook

# --- Version History sheet: log the v1.1 entry for this review pass ---
$wsHistory = $wb.Worksheets.Item("Version History")
$wsHistory.Range("A3").Value = "v1.1"
$wsHistory.Range("B3").Value = "Gehad Ashry"
$wsHistory.Range("C3").Value = "Check reviews"
$wsHistory.Range("D3").Value = 45769

# --- Reviews sheet: update Owner Status column (I) after checking reviews ---
$wsReviews = $wb.Worksheets.Item("LH-TC-SYSTEMCONSTRAINS-REVIEWS")
$wsReviews.Range("I2").Value = "Closed"
$wsReviews.Range("I3").Value = "Closed"
$wsReviews.Range("I4").Value = "Closed"
$wsReviews.Range("I5").Value = "Closed"
$wsReviews.Range("I6").Value = "NotApplicable"
$wsReviews.Range("I7").Value = "Closed"
$wsReviews.Range("I8").Value = "Closed"

# --- Final UI state: land on the Version History sheet with C8 selected ---
$wsReviews.Activate()
$wsReviews.Range("I8").Select()
$wsHistory.Activate()
$wsHistory.Range("C8").Select()
